# Auto-generated COM-interop script to apply the quant Mid Cap Fund holdings update
# (adds a Status column, refreshes Jan/Dec/Oct values + MoM/QoQ, reorders/
#  appends the fully-exited holdings rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new 'Status' column at D, shifting Jan_2026..QoQ one column to the right
$ws.Columns.Item(4).Insert()

# 2) Fix up the header row (A1:I1)
$ws.Cells.Item(1,4).Value = "Status"
$ws.Cells.Item(1,7).Value = "Oct_2025"

# 3) Update data rows 2-21 (same ISIN order as before, new Status column + refreshed month data)
$ws.Cells.Item(2,4).Value = "Adding Consistently"
$ws.Cells.Item(2,5).Value = 9.897539
$ws.Cells.Item(2,6).Value = 8.763819
$ws.Cells.Item(2,7).Value = 7.974271
$ws.Cells.Item(2,8).Value = 1.13372
$ws.Cells.Item(2,9).Value = 1.923268
$ws.Cells.Item(3,4).Value = "Reducing Consistently"
$ws.Cells.Item(3,5).Value = 8.547462
$ws.Cells.Item(3,6).Value = 8.979951
$ws.Cells.Item(3,7).Value = 8.720279
$ws.Cells.Item(3,8).Value = -0.4324890000000003
$ws.Cells.Item(3,9).Value = -0.1728170000000002
$ws.Cells.Item(4,4).Value = "Adding Consistently"
$ws.Cells.Item(4,5).Value = 7.117608
$ws.Cells.Item(4,6).Value = 6.668885
$ws.Cells.Item(4,7).Value = 6.713747
$ws.Cells.Item(4,8).Value = 0.4487229999999993
$ws.Cells.Item(4,9).Value = 0.403861
$ws.Cells.Item(5,4).Value = "Adding Consistently"
$ws.Cells.Item(5,5).Value = 6.003223
$ws.Cells.Item(5,6).Value = 5.987753
$ws.Cells.Item(5,7).Value = 1.862603
$ws.Cells.Item(5,8).Value = 0.01547000000000054
$ws.Cells.Item(5,9).Value = 4.14062
$ws.Cells.Item(6,4).Value = "Reducing Consistently"
$ws.Cells.Item(6,5).Value = 5.903756
$ws.Cells.Item(6,6).Value = 6.242897
$ws.Cells.Item(6,7).Value = 7.106878
$ws.Cells.Item(6,8).Value = -0.3391410000000006
$ws.Cells.Item(6,9).Value = -1.203122
$ws.Cells.Item(7,4).Value = "Reducing Consistently"
$ws.Cells.Item(7,5).Value = 5.563413
$ws.Cells.Item(7,6).Value = 6.01995
$ws.Cells.Item(7,7).Value = 5.618652
$ws.Cells.Item(7,8).Value = -0.456537
$ws.Cells.Item(7,9).Value = -0.05523900000000026
$ws.Cells.Item(8,4).Value = "Adding Consistently"
$ws.Cells.Item(8,5).Value = 4.177577
$ws.Cells.Item(8,6).Value = 4.105208
$ws.Cells.Item(8,7).Value = 1.402193
$ws.Cells.Item(8,8).Value = 0.07236900000000013
$ws.Cells.Item(8,9).Value = 2.775384
$ws.Cells.Item(9,4).Value = "Adding Consistently"
$ws.Cells.Item(9,5).Value = 3.948508
$ws.Cells.Item(9,6).Value = 3.538621
$ws.Cells.Item(9,7).Value = 3.377069
$ws.Cells.Item(9,8).Value = 0.4098869999999999
$ws.Cells.Item(9,9).Value = 0.5714389999999998
$ws.Cells.Item(10,4).Value = "Adding Consistently"
$ws.Cells.Item(10,5).Value = 3.113632
$ws.Cells.Item(10,6).Value = 2.845461
$ws.Cells.Item(10,7).Value = 1.495565
$ws.Cells.Item(10,8).Value = 0.2681710000000002
$ws.Cells.Item(10,9).Value = 1.618067
$ws.Cells.Item(11,4).Value = "Adding Consistently"
$ws.Cells.Item(11,5).Value = 2.346329
$ws.Cells.Item(11,6).Value = 2.336732
$ws.Cells.Item(11,7).Value = 2.247703
$ws.Cells.Item(11,8).Value = 0.009596999999999856
$ws.Cells.Item(11,9).Value = 0.09862599999999988
$ws.Cells.Item(12,4).Value = "Reducing Consistently"
$ws.Cells.Item(12,5).Value = 2.030154
$ws.Cells.Item(12,6).Value = 6.012765
$ws.Cells.Item(12,7).Value = 5.378936
$ws.Cells.Item(12,8).Value = -3.982611
$ws.Cells.Item(12,9).Value = -3.348782
$ws.Cells.Item(13,4).Value = "Reducing"
$ws.Cells.Item(13,5).Value = 1.593499
$ws.Cells.Item(13,6).Value = 1.831049
$ws.Cells.Item(13,7).Value = 0.825861
$ws.Cells.Item(13,8).Value = -0.2375499999999999
$ws.Cells.Item(13,9).Value = 0.767638
$ws.Cells.Item(14,4).Value = "Adding Consistently"
$ws.Cells.Item(14,5).Value = 1.570087
$ws.Cells.Item(14,6).Value = 1.561581
$ws.Cells.Item(14,7).Value = 1.356743
$ws.Cells.Item(14,8).Value = 0.008505999999999903
$ws.Cells.Item(14,9).Value = 0.213344
$ws.Cells.Item(15,4).Value = "Adding"
$ws.Cells.Item(15,5).Value = 1.49927
$ws.Cells.Item(15,6).Value = 1.496639
$ws.Cells.Item(15,7).Value = 1.56779
$ws.Cells.Item(15,8).Value = 0.00263100000000005
$ws.Cells.Item(15,9).Value = -0.06851999999999991
$ws.Cells.Item(16,4).Value = "Adding Consistently"
$ws.Cells.Item(16,5).Value = 1.441495
$ws.Cells.Item(16,6).Value = 0.451678
$ws.Cells.Item(16,7).Value = 0.472883
$ws.Cells.Item(16,8).Value = 0.989817
$ws.Cells.Item(16,9).Value = 0.968612
$ws.Cells.Item(17,4).Value = "Fresh Entry"
$ws.Cells.Item(17,5).Value = 1.333105
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 0
$ws.Cells.Item(17,8).Value = 1.333105
$ws.Cells.Item(17,9).Value = 1.333105
$ws.Cells.Item(18,4).Value = "Reducing Consistently"
$ws.Cells.Item(18,5).Value = 1.181887
$ws.Cells.Item(18,6).Value = 2.43113
$ws.Cells.Item(18,7).Value = 5.001537
$ws.Cells.Item(18,8).Value = -1.249243
$ws.Cells.Item(18,9).Value = -3.81965
$ws.Cells.Item(19,4).Value = "Reducing"
$ws.Cells.Item(19,5).Value = 1.150061
$ws.Cells.Item(19,6).Value = 1.152939
$ws.Cells.Item(19,7).Value = 0.675934
$ws.Cells.Item(19,8).Value = -0.002877999999999936
$ws.Cells.Item(19,9).Value = 0.474127
$ws.Cells.Item(20,4).Value = "Adding Consistently"
$ws.Cells.Item(20,5).Value = 1.118878
$ws.Cells.Item(20,6).Value = 1.048692
$ws.Cells.Item(20,7).Value = 1.037306
$ws.Cells.Item(20,8).Value = 0.07018600000000008
$ws.Cells.Item(20,9).Value = 0.08157199999999998
$ws.Cells.Item(21,4).Value = "Fresh Entry"
$ws.Cells.Item(21,5).Value = 0.717933
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = 0.717933
$ws.Cells.Item(21,9).Value = 0.717933

# 4) Clear the old rows 22-26 (exited holdings) - they get replaced/reordered below
$ws.Range("A22:I26").ClearContents()

# 5) Write rows 22-29: reordered + new 'Complete Exit' holdings
$ws.Cells.Item(22,1).Value = "INE335Y01020"
$ws.Cells.Item(22,2).Value = "Indian Railway Catering & Tourism Corp"
$ws.Cells.Item(22,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(22,4).Value = "Complete Exit"
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 1.395242
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(22,9).Value = -1.395242
$ws.Cells.Item(23,1).Value = "INE376G01013"
$ws.Cells.Item(23,2).Value = "Biocon Ltd"
$ws.Cells.Item(23,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(23,4).Value = "Complete Exit"
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 0.381318
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = -0.381318
$ws.Cells.Item(23,9).Value = 0
$ws.Cells.Item(24,1).Value = "INE245A01021"
$ws.Cells.Item(24,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(24,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(24,4).Value = "Complete Exit"
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 2.202735
$ws.Cells.Item(24,7).Value = 2.220652
$ws.Cells.Item(24,8).Value = -2.202735
$ws.Cells.Item(24,9).Value = -2.220652
$ws.Cells.Item(25,1).Value = "INE154A01025"
$ws.Cells.Item(25,2).Value = "ITC Limited"
$ws.Cells.Item(25,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(25,4).Value = "Complete Exit"
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 2.500818
$ws.Cells.Item(25,7).Value = 2.465386
$ws.Cells.Item(25,8).Value = -2.500818
$ws.Cells.Item(25,9).Value = -2.465386
$ws.Cells.Item(26,1).Value = "INE14LE01019"
$ws.Cells.Item(26,2).Value = "Aditya Birla Lifestyle Brands Limited"
$ws.Cells.Item(26,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(26,4).Value = "Complete Exit"
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 0.433233
$ws.Cells.Item(26,8).Value = 0
$ws.Cells.Item(26,9).Value = -0.433233
$ws.Cells.Item(27,1).Value = "INE115A01026"
$ws.Cells.Item(27,2).Value = "LIC Housing Finance Ltd"
$ws.Cells.Item(27,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(27,4).Value = "Complete Exit"
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 1.475937
$ws.Cells.Item(27,7).Value = 1.476151
$ws.Cells.Item(27,8).Value = -1.475937
$ws.Cells.Item(27,9).Value = -1.476151
$ws.Cells.Item(28,1).Value = "INE094A01015"
$ws.Cells.Item(28,2).Value = "Hindustan Petroleum Corporation Ltd"
$ws.Cells.Item(28,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(28,4).Value = "Complete Exit"
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 0
$ws.Cells.Item(28,7).Value = 3.216334
$ws.Cells.Item(28,8).Value = 0
$ws.Cells.Item(28,9).Value = -3.216334
$ws.Cells.Item(29,1).Value = "INE326A01037"
$ws.Cells.Item(29,2).Value = "Lupin Limited"
$ws.Cells.Item(29,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(29,4).Value = "Complete Exit"
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = 0
$ws.Cells.Item(29,7).Value = 0.992746
$ws.Cells.Item(29,8).Value = 0
$ws.Cells.Item(29,9).Value = -0.992746

